$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.614229559898376
$ws.Range("B1").Value = 4.238815307617188
$ws.Range("C1").Value = 3.576089859008789
$ws.Range("D1").Value = 1.797138094902039
$ws.Range("E1").Value = 1.041212558746338
